$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.827.75"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.836.36"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "308.82"
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  +3.52%  "

$ws.Range("D8").Value = "0.3652"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").Value = "0.07138"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "0.9178"
$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "19.50"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07649"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.781.51"
$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").Value = "5.277"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "6.386"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").Value = "87.86"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "0.000008617"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "26.883.97"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("E21").Value = "  +2.16%  "

$ws.Range("D22").Value = "5.000"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("D24").Value = "1.918"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "151.61"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "18.18"
$ws.Range("E26").Value = "  +2.24%  "

$ws.Range("D27").Value = "2.004"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "114.04"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("D29").Value = "4.870"
$ws.Range("E29").Value = "  +1.01%  "

$ws.Range("D30").Value = "0.08814"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").Value = "3.206"
$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("D32").Value = "1.172"
$ws.Range("E32").Value = "  +5.65%  "

$ws.Range("D33").Value = "0.7430"
$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").Value = "2.742"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").Value = "4.461"
$ws.Range("E35").Value = "  +0.84%  "

$ws.Range("D36").Value = "1.087"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").Value = "0.05205"
$ws.Range("E38").Value = "  +2.45%  "

$ws.Range("D39").Value = "2.953"
$ws.Range("E39").Value = "  +1.39%  "

$ws.Range("D40").Value = "0.5175"
$ws.Range("E40").Value = "  +1.85%  "

$ws.Range("D41").Value = "6.952"
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").Value = "0.1509"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "8.130"
$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("D44").Value = "10.47"
$ws.Range("E44").Value = "  +5.23%  "

$ws.Range("D45").Value = "0.4688"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "101.76"
$ws.Range("E47").Value = "  +2.85%  "

$ws.Range("D48").Value = "1.589"
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("D49").Value = "64.78"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("D50").Value = "0.06030"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").Value = "0.8833"
$ws.Range("E51").Value = "  +4.50%  "
